# Update the "Förändrad" (C) date column for all data rows (2-150) from
# 45184 to 45186, and add a friendly-name second argument (the
# "Beteckning" in column A) to the HYPERLINK() formulas in columns
# S, T, V, W, X, Y for the rows that have them (2-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (as worksheet column numbers) that contain HYPERLINK
# formulas which need a friendly-name second argument added.
# A=1, B=2, C=3, ... S=19, T=20, U=21, V=22, W=23, X=24, Y=25
$hyperlinkCols = @(19, 20, 22, 23, 24, 25)

for ($r = 2; $r -le 150; $r++) {

    # 1) Update the "Förändrad" date in column C (3) for every data row.
    $ws.Cells.Item($r, 3).Value = 45186

    # 2) For the rows that have HYPERLINK formulas (rows 2-16), add the
    #    Beteckning (column A) as the friendly-name second argument.
    if ($r -le 16) {
        $beteckning = $ws.Cells.Item($r, 1).Value2

        foreach ($c in $hyperlinkCols) {
            $cell = $ws.Cells.Item($r, $c)
            $formula = $cell.Formula
            if ($formula) {
                if ($formula -match '^(=HYPERLINK\(".*?")\)$') {
                    $newFormula = $matches[1] + ', "' + $beteckning + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
